{"js": "// Update the worksheet date and every \"three-digit \u00f7 one-digit\" answer\n// cell to the new values from the commit, e.g.\n//   \"2025-02-25 Tuesday\"  -> \"2025-02-26 Wednesday\"\n//   \"312\u00f74=78, 0\"         -> \"702\u00f77=100, 2\"\n//   ... (one pair per changed <w:t> run; every \"old\" string is unique in\n//   the document, so body.search + insertText(\"Replace\") is safe).\nconst pairs = [\n  [\"2025-02-25 Tuesday\", \"2025-02-26 Wednesday\"],\n  [\"312\u00f74=78, 0\", \"702\u00f77=100, 2\"],\n  [\"130\u00f78=16, 2\", \"808\u00f74=202, 0\"],\n  [\"374\u00f75=74, 4\", \"333\u00f79=37, 0\"],\n  [\"408\u00f72=204, 0\", \"821\u00f79=91, 2\"],\n  [\"543\u00f75=108, 3\", \"379\u00f78=47, 3\"],\n  [\"178\u00f77=25, 3\", \"959\u00f72=479, 1\"],\n  [\"411\u00f78=51, 3\", \"396\u00f72=198, 0\"],\n  [\"516\u00f76=86, 0\", \"314\u00f72=157, 0\"],\n  [\"914\u00f78=114, 2\", \"731\u00f73=243, 2\"],\n  [\"558\u00f77=79, 5\", \"443\u00f77=63, 2\"],\n  [\"909\u00f79=101, 0\", \"629\u00f77=89, 6\"],\n  [\"825\u00f74=206, 1\", \"900\u00f76=150, 0\"],\n  [\"214\u00f76=35, 4\", \"702\u00f78=87, 6\"],\n  [\"690\u00f72=345, 0\", \"509\u00f75=101, 4\"],\n  [\"669\u00f72=334, 1\", \"582\u00f76=97, 0\"],\n  [\"294\u00f75=58, 4\", \"699\u00f72=349, 1\"],\n  [\"442\u00f75=88, 2\", \"242\u00f76=40, 2\"],\n  [\"762\u00f74=190, 2\", \"810\u00f76=135, 0\"],\n  [\"334\u00f73=111, 1\", \"872\u00f79=96, 8\"],\n  [\"330\u00f79=36, 6\", \"165\u00f75=33, 0\"],\n  [\"709\u00f79=78, 7\", \"283\u00f72=141, 1\"],\n  [\"520\u00f75=104, 0\", \"668\u00f75=133, 3\"],\n  [\"142\u00f79=15, 7\", \"188\u00f72=94, 0\"],\n  [\"253\u00f73=84, 1\", \"684\u00f72=342, 0\"],\n  [\"551\u00f76=91, 5\", \"472\u00f79=52, 4\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace every match (expected to be exactly one per pair in this doc).\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"three-digit \u00f7 one-digit\" answer\n# cell to the new values from the commit, e.g.\n#   \"2025-02-25 Tuesday\"  -> \"2025-02-26 Wednesday\"\n#   \"312\u00f74=78, 0\"         -> \"702\u00f77=100, 2\"\n#   ... (one pair per changed run; every \"old\" string is unique in the\n#   document, so a simple Find/Replace per pair is safe).\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-02-25 Tuesday\", \"2025-02-26 Wednesday\"),\n  @(\"312\u00f74=78, 0\", \"702\u00f77=100, 2\"),\n  @(\"130\u00f78=16, 2\", \"808\u00f74=202, 0\"),\n  @(\"374\u00f75=74, 4\", \"333\u00f79=37, 0\"),\n  @(\"408\u00f72=204, 0\", \"821\u00f79=91, 2\"),\n  @(\"543\u00f75=108, 3\", \"379\u00f78=47, 3\"),\n  @(\"178\u00f77=25, 3\", \"959\u00f72=479, 1\"),\n  @(\"411\u00f78=51, 3\", \"396\u00f72=198, 0\"),\n  @(\"516\u00f76=86, 0\", \"314\u00f72=157, 0\"),\n  @(\"914\u00f78=114, 2\", \"731\u00f73=243, 2\"),\n  @(\"558\u00f77=79, 5\", \"443\u00f77=63, 2\"),\n  @(\"909\u00f79=101, 0\", \"629\u00f77=89, 6\"),\n  @(\"825\u00f74=206, 1\", \"900\u00f76=150, 0\"),\n  @(\"214\u00f76=35, 4\", \"702\u00f78=87, 6\"),\n  @(\"690\u00f72=345, 0\", \"509\u00f75=101, 4\"),\n  @(\"669\u00f72=334, 1\", \"582\u00f76=97, 0\"),\n  @(\"294\u00f75=58, 4\", \"699\u00f72=349, 1\"),\n  @(\"442\u00f75=88, 2\", \"242\u00f76=40, 2\"),\n  @(\"762\u00f74=190, 2\", \"810\u00f76=135, 0\"),\n  @(\"334\u00f73=111, 1\", \"872\u00f79=96, 8\"),\n  @(\"330\u00f79=36, 6\", \"165\u00f75=33, 0\"),\n  @(\"709\u00f79=78, 7\", \"283\u00f72=141, 1\"),\n  @(\"520\u00f75=104, 0\", \"668\u00f75=133, 3\"),\n  @(\"142\u00f79=15, 7\", \"188\u00f72=94, 0\"),\n  @(\"253\u00f73=84, 1\", \"684\u00f72=342, 0\"),\n  @(\"551\u00f76=91, 5\", \"472\u00f79=52, 4\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n  #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap:=wdFindContinue,\n  #   Format, ReplaceWith, Replace:=wdReplaceAll)\n  $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Could not find text: $oldText\"\n  }\n}\n\n"}
